$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell (untouched, style s="3", General format) used to restore
# formatting on cells whose new text value looks like a percentage, which
# Excel would otherwise auto-convert into a numeric percent value.
$fmtRef = $ws.Range("H2")

$ws.Range("E2").Value = "2026-02-07 07:17:43"
$ws.Range("E3").Value = "2026-02-07 07:17:45"
$ws.Range("N3").Value = "-7.6 °C 6:38 TU"
$ws.Range("O3").Value = "-6.0 °C"
$ws.Range("E4").Value = "2026-02-07 07:17:48"
$ws.Range("J4").Value = "1001.5 hPa"
$ws.Range("N4").Value = "9.3 °C 6:59 TU"
$ws.Range("O4").Value = "11.1 °C"
$ws.Range("E5").Value = "2026-02-07 07:17:50"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "75%"
$fmtRef.Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("J5").Value = "1001.5 hPa"
$ws.Range("O5").Value = "8.3 °C"
$ws.Range("E6").Value = "2026-02-07 07:17:53"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "57%"
$fmtRef.Copy() | Out-Null
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("J6").Value = "1003.1 hPa"
$ws.Range("E7").Value = "2026-02-07 07:17:55"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "74%"
$fmtRef.Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("J7").Value = "1002.8 hPa"
$ws.Range("L7").Value = "50.8 km/h - 345º 6:54 TU"
$ws.Range("E8").Value = "2026-02-07 07:17:57"
$ws.Range("O8").Value = "3.9 °C"
$ws.Range("E9").Value = "2026-02-07 07:17:59"
$ws.Range("O9").Value = "1.3 °C"
$ws.Range("E10").Value = "2026-02-07 07:18:02"
$ws.Range("E11").Value = "2026-02-07 07:18:04"
$ws.Range("E12").Value = "2026-02-07 07:18:07"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "69%"
$fmtRef.Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("E13").Value = "2026-02-07 07:18:09"
$ws.Range("E14").Value = "2026-02-07 07:18:11"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "72%"
$fmtRef.Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("O14").Value = "-5.9 °C"
$ws.Range("E15").Value = "2026-02-07 07:18:13"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "89%"
$fmtRef.Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("J15").Value = "1001.8 hPa"
$ws.Range("O15").Value = "5.7 °C"
$ws.Range("E16").Value = "2026-02-07 07:18:16"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "92%"
$fmtRef.Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("E17").Value = "2026-02-07 07:18:18"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "99%"
$fmtRef.Copy() | Out-Null
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("I17").Value = "0.2 mm"
$ws.Range("J17").Value = "1005.1 hPa"
$ws.Range("N17").Value = "2.2 °C 6:58 TU"
$ws.Range("E18").Value = "2026-02-07 07:18:20"
$ws.Range("O18").Value = "-7.8 °C"
$ws.Range("E19").Value = "2026-02-07 07:18:23"
$ws.Range("J19").Value = "1006.4 hPa"
$ws.Range("N19").Value = "0.3 °C 6:48 TU"
$ws.Range("O19").Value = "3.9 °C"
$ws.Range("E20").Value = "2026-02-07 07:18:25"
$ws.Range("N20").Value = "-5.8 °C 6:53 TU"
$ws.Range("O20").Value = "-4.8 °C"
$ws.Range("E21").Value = "2026-02-07 07:18:28"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "78%"
$fmtRef.Copy() | Out-Null
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("J21").Value = "1002.2 hPa"
$ws.Range("N21").Value = "2.0 °C 6:46 TU"
$ws.Range("O21").Value = "5.8 °C"
$ws.Range("E22").Value = "2026-02-07 07:18:30"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "89%"
$fmtRef.Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("L22").Value = "17.6 km/h - 359º 6:57 TU"
$ws.Range("M22").Value = "9.6 °C 6:59 TU"
$ws.Range("O22").Value = "6.4 °C"
$ws.Range("E23").Value = "2026-02-07 07:18:32"
$ws.Range("J23").Value = "1001.6 hPa"
$ws.Range("N23").Value = "5.8 °C 6:46 TU"
$ws.Range("O23").Value = "7.4 °C"
$ws.Range("E24").Value = "2026-02-07 07:18:35"
$ws.Range("J24").Value = "1000.9 hPa"
$ws.Range("N24").Value = "9.5 °C 6:42 TU"
$ws.Range("O24").Value = "10.2 °C"
$ws.Range("E25").Value = "2026-02-07 07:18:37"
$ws.Range("J25").Value = "1005.4 hPa"
$ws.Range("E26").Value = "2026-02-07 07:18:40"
$ws.Range("O26").Value = "-2.5 °C"
$ws.Range("E27").Value = "2026-02-07 07:18:42"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "93%"
$fmtRef.Copy() | Out-Null
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("J27").Value = "1001.4 hPa"
$ws.Range("L27").Value = "21.2 km/h - 27º 6:45 TU"
$ws.Range("M27").Value = "11.5 °C 6:47 TU"
$ws.Range("O27").Value = "8.9 °C"
$ws.Range("E28").Value = "2026-02-07 07:18:45"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "90%"
$fmtRef.Copy() | Out-Null
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("J28").Value = "1004.3 hPa"
$ws.Range("N28").Value = "1.6 °C 6:58 TU"
$ws.Range("O28").Value = "3.0 °C"
$ws.Range("E29").Value = "2026-02-07 07:18:47"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "59%"
$fmtRef.Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("N29").Value = "8.3 °C 6:50 TU"
$ws.Range("O29").Value = "10.7 °C"
$ws.Range("E30").Value = "2026-02-07 07:18:49"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "81%"
$fmtRef.Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("E31").Value = "2026-02-07 07:18:52"
$ws.Range("J31").Value = "1006.0 hPa"
$ws.Range("O31").Value = "3.5 °C"
$ws.Range("E32").Value = "2026-02-07 07:18:54"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "57%"
$fmtRef.Copy() | Out-Null
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("J32").Value = "1004.6 hPa"
$ws.Range("E33").Value = "2026-02-07 07:18:56"
$ws.Range("E34").Value = "2026-02-07 07:18:59"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "79%"
$fmtRef.Copy() | Out-Null
$ws.Range("H34").PasteSpecial(-4122)
$ws.Range("N34").Value = "3.8 °C 6:43 TU"
$ws.Range("O34").Value = "6.0 °C"
$ws.Range("E35").Value = "2026-02-07 07:19:01"
$ws.Range("O35").Value = "-6.1 °C"
$ws.Range("E36").Value = "2026-02-07 07:19:04"
$ws.Range("N36").Value = "2.9 °C 6:55 TU"
$ws.Range("O36").Value = "4.5 °C"

$excel.CutCopyMode = $false

